# Landscaping Data.xlsx - "Add files via upload" edit
# Adds 14 new observation rows (128-141, date 5/28/2025) to Sheet1,
# extends the ABS(Low-High) helper formula in column F down through the
# new rows, and updates the saved view (scroll position / selection /
# window geometry) to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Seed rows 128:141 from row 127 so the new rows inherit the same
#        cell formatting (in particular the date number format on column A)
#        instead of picking up the engine's generic default style. ---
$ws.Range("A127:T127").Copy($ws.Range("A128:T141"))

# --- 2. New row data (Date, Plant_Type, Plant_Size, Low, High, [Temp_Diff
#        formula], Rain, Growth, Pruned, Quadrant, Shade, UV, Humidity,
#        Dew_Point, Pressure, Wind_Gust, Cloud_Cover, Visibility, AQI,
#        Pollen). Mirrors the two new daily observation blocks (date
#        45805 / 5-28-2025) appended to the tracking sheet. ---
$rows = @(
    @(45805,"Flowering","Large",56,58,1.05,0.1,"No",2,"Dark",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Nonflowering","Medium",56,58,1.05,0.2,"No",3,"Neutral",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Nonflowering","Small",56,58,1.05,0.3,"No",3,"Neutral",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Nonflowering","Medium",56,58,1.05,0.2,"No",3,"Dark",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Nonflowering","Medium",56,58,1.05,0.5,"No",3,"Neutral",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Nonflowering","Large",56,58,1.05,0.2,"No",4,"Dark",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Tree","Medium",56,58,1.05,0.15,"No",1,"Dark",3,0.93,56,30.07,18,0.95,8.7,26,39),
    @(45805,"Flowering","Large",56,71,0.27,0.2,"No",2,"Dark",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Nonflowering","Medium",56,71,0.27,0.5,"No",3,"Bright",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Nonflowering","Small",56,71,0.27,0.13,"No",3,"Bright",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Nonflowering","Medium",56,71,0.27,0.6,"No",3,"Neutral",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Nonflowering","Medium",56,71,0.27,1.25,"No",3,"Neutral",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Nonflowering","Large",56,71,0.27,0.8,"No",4,"Dark",3,0.68,60,29.99,15,0.83,9.9,35,33),
    @(45805,"Tree","Medium",56,71,0.27,$null,"No",1,"Neutral",3,0.68,60,29.99,15,0.83,9.9,35,33)
)

$r = 128
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]                 # A Date
    $ws.Cells.Item($r, 2).Value = $row[1]                 # B Plant_Type
    $ws.Cells.Item($r, 3).Value = $row[2]                 # C Plant_Size
    $ws.Cells.Item($r, 4).Value = $row[3]                 # D Low
    $ws.Cells.Item($r, 5).Value = $row[4]                 # E High
    $ws.Cells.Item($r, 6).Formula = "=ABS(D$r-E$r)"       # F Temp_Diff
    $ws.Cells.Item($r, 7).Value = $row[5]                 # G Rain
    if ($r -eq 141) {
        $ws.Cells.Item($r, 8).Formula = "=4/3"            # H Growth (formula on last row)
    } else {
        $ws.Cells.Item($r, 8).Value = $row[6]             # H Growth
    }
    $ws.Cells.Item($r, 9).Value = $row[7]                 # I Pruned
    $ws.Cells.Item($r, 10).Value = $row[8]                # J Quadrant
    $ws.Cells.Item($r, 11).Value = $row[9]                # K Shade
    $ws.Cells.Item($r, 12).Value = $row[10]               # L UV
    $ws.Cells.Item($r, 13).Value = $row[11]               # M Humidity
    $ws.Cells.Item($r, 14).Value = $row[12]               # N Dew_Point
    $ws.Cells.Item($r, 15).Value = $row[13]               # O Pressure
    $ws.Cells.Item($r, 16).Value = $row[14]               # P Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row[15]               # Q Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row[16]               # R Visibility
    $ws.Cells.Item($r, 19).Value = $row[17]               # S AQI
    $ws.Cells.Item($r, 20).Value = $row[18]               # T Pollen
    $r++
}

# --- 3. Update the saved scroll position / selection to match where the
#        author ended up after pasting the new data. ---
$ws.Application.Goto($ws.Range("A133"))
$ws.Range("Q135:Q141").Select()

# --- 4. Update the workbook window geometry saved with the file. ---
$excel.Left = 1236
$excel.Top = 3048
$excel.Width = 21156
$excel.Height = 7128
